# Populate row 1 of the (previously empty) sheet with the class-schedule
# record described by the commit. Columns are written in sheet order;
# numeric-looking columns get real numbers, everything else is literal text
# (Excel's text-to-number/date autoconversion is avoided by keeping the
# look-alike values, like dates and dashes, as plain strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 257
$ws.Range("B1").Value = "DEI"
$ws.Range("C1").Value = "Redes de Computadoras II"
$ws.Range("D1").Value = "---"
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = "IIN"
$ws.Range("G1").Value = "-- --"
$ws.Range("H1").Value = 2008
$ws.Range("I1").Value = "T"
$ws.Range("J1").Value = "TQ"
$ws.Range("K1").Value = "EDUCA"
$ws.Range("L1").Value = "Ms."
$ws.Range("M1").Value = "Amarilla Fleitas"
$ws.Range("N1").Value = "Gustavo Osman"
$ws.Range("O1").Value = "gamarilla@pol.una.py"

$ws.Range("P1").Value = "Mie 03/09/25"
$ws.Range("Q1").Value = 0.75
$ws.Range("Q1").NumberFormat = "h:mm:ss"

$ws.Range("S1").Value = "Mie 29/10/25"
$ws.Range("T1").Value = 0.75
$ws.Range("T1").NumberFormat = "h:mm:ss"

$ws.Range("V1").Value = "Mie 19/11/25"
$ws.Range("W1").Value = 0.75
$ws.Range("W1").NumberFormat = "h:mm:ss"

$ws.Range("AA1").Value = "Mie 03/12/25"
$ws.Range("AB1").Value = 0.75
$ws.Range("AB1").NumberFormat = "h:mm:ss"

$ws.Range("AF1").Value = "Ms. Gustavo Osman Amarilla Fleitas"
$ws.Range("AG1").Value = "Ms. Claudio Nil Barúa Acosta"
$ws.Range("AH1").Value = "Lic. María Luisa Guanes Romero"

$ws.Range("AO1").Value = "I02"
$ws.Range("AP1").Value = "14:15 - 16:30"
$ws.Range("AS1").Value = "I02"
$ws.Range("AT1").Value = "07:30 - 10:30"

# The saved sheet's used range stretches one column further right (through
# AU) than the last populated cell (AT). Nudge the sheet's extent out to AU1
# (formatting-only touch, no value) so the exported dimension matches.
$ws.Range("AU1").Font.Bold = $false
